# Add data for 2021-10-06
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-09-28"

# Update the September header label
$ws.Range("A10").Value = "September (through 09-28)"

# Update September row (row 10) values
$ws.Range("C10").Value = 42
$ws.Range("D10").Value = 72
$ws.Range("E10").Value = 52
$ws.Range("F10").Value = 68
$ws.Range("G10").Value = 107
$ws.Range("H10").Value = 168

# Update Total row (row 11) values
$ws.Range("C11").Value = 423
$ws.Range("D11").Value = 623
$ws.Range("E11").Value = 542
$ws.Range("F11").Value = 417
$ws.Range("G11").Value = 891
$ws.Range("H11").Value = 1238
